$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.166.90"
$ws.Range("E2").Value = "  +1.17%  "
$ws.Range("D3").Value = "1.640.59"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("E6").Value = "  +2.54%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0625"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.94"
$ws.Range("D10").Style = "Normal"
$ws.Range("E11").Value = "  +0.42%  "
$ws.Range("D12").Value = "1.871.18"
$ws.Range("D13").Value = "1.653.86"
$ws.Range("E13").Value = "  +1.10%  "
$ws.Range("E14").Value = "  +0.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.94"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.51%  "
$ws.Range("D17").Value = "27.167.18"
$ws.Range("E17").Value = "  +1.15%  "
$ws.Range("D18").Value = "0.0₃0739"
$ws.Range("E18").Value = "  +1.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "216.88"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.27%  "
$ws.Range("E20").Value = "  +0.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.35%  "
$ws.Range("E23").Value = "  +2.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.04"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.43"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.96%  "
$ws.Range("E28").Value = "  +0.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.67"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.80%  "
$ws.Range("E30").Value = "  +0.74%  "
$ws.Range("E31").Value = "  +0.24%  "
$ws.Range("E32").Value = "  +1.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.10%  "
$ws.Range("E34").Value = "  +0.85%  "
$ws.Range("D35").Value = "1.297.16"
$ws.Range("E35").Value = "  +3.38%  "
$ws.Range("E36").Value = "  +1.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0176"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.550"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.05%  "
$ws.Range("E39").Value = "  +2.91%  "
$ws.Range("E40").Value = "  +0.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.809"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.22"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.70%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.30"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.83%  "
$ws.Range("D44").Value = "1.780.97"
$ws.Range("E44").Value = "  +0.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.66"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.76"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.59"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.00%  "
$ws.Range("E48").Value = "  +1.82%  "
$ws.Range("E49").Value = "  -0.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.62"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0962"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.16%  "
